# "updated main GSC export data"
# Append the next day's row (2025-12-31) to the bottom of the "Chart"
# sheet's data table, carrying forward the same Clicks/Impressions
# totals as the prior day (0 / 28), exactly like the GSC export script
# that produces this workbook normally does.

$wb    = $excel.ActiveWorkbook
$chart = $wb.Worksheets.Item("Chart")

$newRow = 87

# Write the date as text (not an auto-converted date serial) so it
# matches the existing column A cells, which are all plain strings.
$chart.Cells.Item($newRow, 1).NumberFormat = "@"
$chart.Cells.Item($newRow, 1).Value = "2025-12-31"

# Re-apply the (default) formatting of the row above so the new cell
# doesn't keep a stray "text" number format applied above, matching
# the unformatted style used throughout the rest of the table.
$chart.Cells.Item($newRow - 1, 1).Copy() | Out-Null
$chart.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null

$chart.Cells.Item($newRow, 2).Value = 0
$chart.Cells.Item($newRow, 3).Value = 28
